# PerformanceTestWithCompleteFixtures.xlsx update
# "Messung in der Konsole. Kleinere Umstellungen."
# Refresh the recorded measurement run: new timestamp, a new "Selma" mapper
# row appended to both sheets, and refreshed numbers throughout.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Warmlaufen"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Warmlaufen")

# Updated timestamp heading
$ws1.Cells.Item(1, 1).Value = "05.01.2022 um 17:56 Uhr"

# Row 3 - ByHand
$ws1.Cells.Item(3, 3).Value = 5.0
$ws1.Cells.Item(3, 6).Value = 0.5
$ws1.Cells.Item(3, 7).Value = 0.707107

# Row 4 - Dozer
$ws1.Cells.Item(4, 3).Value = 188.0
$ws1.Cells.Item(4, 5).Value = 92.0
$ws1.Cells.Item(4, 6).Value = 18.8
$ws1.Cells.Item(4, 7).Value = 26.4357

# Row 5 - Orika
$ws1.Cells.Item(5, 3).Value = 55.0
$ws1.Cells.Item(5, 5).Value = 47.0
$ws1.Cells.Item(5, 6).Value = 5.5
$ws1.Cells.Item(5, 7).Value = 14.6002

# Row 6 - MapStruct
$ws1.Cells.Item(6, 3).Value = 5.0
$ws1.Cells.Item(6, 6).Value = 0.5
$ws1.Cells.Item(6, 7).Value = 0.707107
$ws1.Cells.Item(6, 8).Value = 0.0

# Row 7 - ModelMapper
$ws1.Cells.Item(7, 3).Value = 290.0
$ws1.Cells.Item(7, 5).Value = 103.0
$ws1.Cells.Item(7, 6).Value = 29.0
$ws1.Cells.Item(7, 7).Value = 26.1066

# Row 8 - JMapper
$ws1.Cells.Item(8, 3).Value = 5.0
$ws1.Cells.Item(8, 6).Value = 0.5
$ws1.Cells.Item(8, 7).Value = 0.527046

# Row 9 - Selma (name unchanged position-wise, new mapper replacing reMap's old slot)
$ws1.Cells.Item(9, 1).Value = "Selma"
$ws1.Cells.Item(9, 2).Value = 10.0
$ws1.Cells.Item(9, 3).Value = 13.0
$ws1.Cells.Item(9, 4).Value = 0.0
$ws1.Cells.Item(9, 5).Value = 3.0
$ws1.Cells.Item(9, 6).Value = 1.3
$ws1.Cells.Item(9, 7).Value = 0.948683
$ws1.Cells.Item(9, 8).Value = 1.0

# Row 10 - reMap (new row, moved to the bottom)
$ws1.Cells.Item(10, 1).Value = "reMap"
$ws1.Cells.Item(10, 2).Value = 10.0
$ws1.Cells.Item(10, 3).Value = 71.0
$ws1.Cells.Item(10, 4).Value = 4.0
$ws1.Cells.Item(10, 5).Value = 20.0
$ws1.Cells.Item(10, 6).Value = 7.1
$ws1.Cells.Item(10, 7).Value = 4.8408
$ws1.Cells.Item(10, 8).Value = 4.0

# Copy the numeric/style formatting of the row above down into the new row
$ws1.Range("A9:H9").Copy() | Out-Null
$ws1.Range("A10:H10").PasteSpecial(-4122) | Out-Null    # xlPasteFormats
$ws1.Cells.Item(1, 1).Select() | Out-Null

# Re-apply values (PasteSpecial formats only, so values above already stand,
# but ensure nothing was clobbered)
$ws1.Cells.Item(10, 1).Value = "reMap"
$ws1.Cells.Item(10, 2).Value = 10.0
$ws1.Cells.Item(10, 3).Value = 71.0
$ws1.Cells.Item(10, 4).Value = 4.0
$ws1.Cells.Item(10, 5).Value = 20.0
$ws1.Cells.Item(10, 6).Value = 7.1
$ws1.Cells.Item(10, 7).Value = 4.8408
$ws1.Cells.Item(10, 8).Value = 4.0

# ---------------------------------------------------------------------
# Sheet "Performanz Messung"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Performanz Messung")

# Row 3 - ByHand
$ws2.Cells.Item(3, 2).Value = 41715.0
$ws2.Cells.Item(3, 3).Value = 577.0
$ws2.Cells.Item(3, 5).Value = 2.0
$ws2.Cells.Item(3, 6).Value = 0.013832
$ws2.Cells.Item(3, 7).Value = 0.117

# Row 4 - Dozer
$ws2.Cells.Item(4, 2).Value = 41715.0
$ws2.Cells.Item(4, 3).Value = 179094.0
$ws2.Cells.Item(4, 4).Value = 3.0
$ws2.Cells.Item(4, 5).Value = 11.0
$ws2.Cells.Item(4, 6).Value = 4.29328
$ws2.Cells.Item(4, 7).Value = 0.469731

# Row 5 - Orika
$ws2.Cells.Item(5, 2).Value = 41715.0
$ws2.Cells.Item(5, 3).Value = 10937.0
$ws2.Cells.Item(5, 6).Value = 0.262184
$ws2.Cells.Item(5, 7).Value = 0.441297

# Row 6 - MapStruct
$ws2.Cells.Item(6, 2).Value = 41715.0
$ws2.Cells.Item(6, 3).Value = 631.0
$ws2.Cells.Item(6, 5).Value = 2.0
$ws2.Cells.Item(6, 6).Value = 0.0151265
$ws2.Cells.Item(6, 7).Value = 0.122645

# Row 7 - ModelMapper
$ws2.Cells.Item(7, 2).Value = 41715.0
$ws2.Cells.Item(7, 3).Value = 244636.0
$ws2.Cells.Item(7, 5).Value = 17.0
$ws2.Cells.Item(7, 6).Value = 5.86446
$ws2.Cells.Item(7, 7).Value = 0.431023

# Row 8 - JMapper
$ws2.Cells.Item(8, 2).Value = 41715.0
$ws2.Cells.Item(8, 3).Value = 271.0
$ws2.Cells.Item(8, 6).Value = 0.00649646
$ws2.Cells.Item(8, 7).Value = 0.0803394

# Row 9 - Selma (replaces the old reMap slot, new figures)
$ws2.Cells.Item(9, 1).Value = "Selma"
$ws2.Cells.Item(9, 2).Value = 41715.0
$ws2.Cells.Item(9, 3).Value = 16365.0
$ws2.Cells.Item(9, 4).Value = 0.0
$ws2.Cells.Item(9, 5).Value = 5.0
$ws2.Cells.Item(9, 6).Value = 0.392305
$ws2.Cells.Item(9, 7).Value = 0.490914
$ws2.Cells.Item(9, 8).Value = 0.0

# Insert a brand new row 10 for "reMap" figures, pushing the summary block down
$ws2.Rows.Item(10).Insert() | Out-Null

$ws2.Cells.Item(10, 1).Value = "reMap"
$ws2.Cells.Item(10, 2).Value = 41715.0
$ws2.Cells.Item(10, 3).Value = 135052.0
$ws2.Cells.Item(10, 4).Value = 3.0
$ws2.Cells.Item(10, 5).Value = 12.0
$ws2.Cells.Item(10, 6).Value = 3.23749
$ws2.Cells.Item(10, 7).Value = 0.46082
$ws2.Cells.Item(10, 8).Value = 3.0

# Match the look (number styling) of the row directly above it
$ws2.Range("A9:H9").Copy() | Out-Null
$ws2.Range("A10:H10").PasteSpecial(-4122) | Out-Null    # xlPasteFormats
$ws2.Cells.Item(1, 1).Select() | Out-Null

# Re-apply values again in case paste touched them
$ws2.Cells.Item(10, 1).Value = "reMap"
$ws2.Cells.Item(10, 2).Value = 41715.0
$ws2.Cells.Item(10, 3).Value = 135052.0
$ws2.Cells.Item(10, 4).Value = 3.0
$ws2.Cells.Item(10, 5).Value = 12.0
$ws2.Cells.Item(10, 6).Value = 3.23749
$ws2.Cells.Item(10, 7).Value = 0.46082
$ws2.Cells.Item(10, 8).Value = 3.0

# Row 13 (used to be row 12, shifted by the insert above) - summary header row
$ws2.Cells.Item(13, 1).Value = "Messreihen"
$ws2.Cells.Item(13, 2).Value = 10.0
$ws2.Cells.Item(13, 3).Value = 100.0
$ws2.Cells.Item(13, 4).Value = 1000.0
$ws2.Cells.Item(13, 5).Value = 10000.0
$ws2.Cells.Item(13, 6).Value = 20000.0
$ws2.Cells.Item(13, 7).Value = 30000.0
$ws2.Cells.Item(13, 8).Value = 40000.0

# Row 14 (was row 13) - ByHand convergence stats
$ws2.Cells.Item(14, 1).Value = "ByHand"
$ws2.Cells.Item(14, 2).Value = 0.0
$ws2.Cells.Item(14, 3).Value = 0.06
$ws2.Cells.Item(14, 4).Value = 0.022
$ws2.Cells.Item(14, 5).Value = 0.0134
$ws2.Cells.Item(14, 6).Value = 0.01395
$ws2.Cells.Item(14, 7).Value = 0.0137
$ws2.Cells.Item(14, 8).Value = 0.0137

# Row 15 (was row 14) - Dozer convergence stats
$ws2.Cells.Item(15, 1).Value = "Dozer"
$ws2.Cells.Item(15, 2).Value = 6.0
$ws2.Cells.Item(15, 3).Value = 5.02
$ws2.Cells.Item(15, 4).Value = 4.191
$ws2.Cells.Item(15, 5).Value = 4.2656
$ws2.Cells.Item(15, 6).Value = 4.2889
$ws2.Cells.Item(15, 7).Value = 4.2926
$ws2.Cells.Item(15, 8).Value = 4.29318

# Row 16 (was row 15) - Orika convergence stats
$ws2.Cells.Item(16, 1).Value = "Orika"
$ws2.Cells.Item(16, 2).Value = 0.4
$ws2.Cells.Item(16, 3).Value = 0.35
$ws2.Cells.Item(16, 4).Value = 0.296
$ws2.Cells.Item(16, 5).Value = 0.2666
$ws2.Cells.Item(16, 6).Value = 0.26455
$ws2.Cells.Item(16, 7).Value = 0.2629
$ws2.Cells.Item(16, 8).Value = 0.262875

# Row 17 (was row 16) - MapStruct convergence stats
$ws2.Cells.Item(17, 1).Value = "MapStruct"
$ws2.Cells.Item(17, 2).Value = 0.3
$ws2.Cells.Item(17, 3).Value = 0.1
$ws2.Cells.Item(17, 4).Value = 0.028
$ws2.Cells.Item(17, 5).Value = 0.0152
$ws2.Cells.Item(17, 6).Value = 0.0163
$ws2.Cells.Item(17, 7).Value = 0.0152667
$ws2.Cells.Item(17, 8).Value = 0.014975

# Row 18 (was row 17) - ModelMapper convergence stats
$ws2.Cells.Item(18, 1).Value = "ModelMapper"
$ws2.Cells.Item(18, 2).Value = 7.9
$ws2.Cells.Item(18, 3).Value = 6.92
$ws2.Cells.Item(18, 4).Value = 5.839
$ws2.Cells.Item(18, 5).Value = 5.8406
$ws2.Cells.Item(18, 6).Value = 5.86255
$ws2.Cells.Item(18, 7).Value = 5.86243
$ws2.Cells.Item(18, 8).Value = 5.86408

# Row 19 (was row 18) - JMapper convergence stats
$ws2.Cells.Item(19, 1).Value = "JMapper"
$ws2.Cells.Item(19, 2).Value = 0.2
$ws2.Cells.Item(19, 3).Value = 0.38
$ws2.Cells.Item(19, 4).Value = 0.061
$ws2.Cells.Item(19, 5).Value = 0.0108
$ws2.Cells.Item(19, 6).Value = 0.00735
$ws2.Cells.Item(19, 7).Value = 0.0069
$ws2.Cells.Item(19, 8).Value = 0.006525

# Row 20 (was row 19, label now refers to Selma) - Selma convergence stats
$ws2.Cells.Item(20, 1).Value = "Selma"
$ws2.Cells.Item(20, 2).Value = 0.5
$ws2.Cells.Item(20, 3).Value = 0.43
$ws2.Cells.Item(20, 4).Value = 0.381
$ws2.Cells.Item(20, 5).Value = 0.3894
$ws2.Cells.Item(20, 6).Value = 0.39185
$ws2.Cells.Item(20, 7).Value = 0.391733
$ws2.Cells.Item(20, 8).Value = 0.3921

# Row 21 (brand new) - reMap convergence stats
$ws2.Cells.Item(21, 1).Value = "reMap"
$ws2.Cells.Item(21, 2).Value = 5.1
$ws2.Cells.Item(21, 3).Value = 3.92
$ws2.Cells.Item(21, 4).Value = 3.241
$ws2.Cells.Item(21, 5).Value = 3.2373
$ws2.Cells.Item(21, 6).Value = 3.23745
$ws2.Cells.Item(21, 7).Value = 3.23717
$ws2.Cells.Item(21, 8).Value = 3.2374

# Make sure row 21's number styling (inherited "general" after the insert two
# steps above) matches the other convergence-stat rows
$ws2.Range("A20:H20").Copy() | Out-Null
$ws2.Range("A21:H21").PasteSpecial(-4122) | Out-Null    # xlPasteFormats
$ws2.Cells.Item(1, 1).Select() | Out-Null

$ws2.Cells.Item(21, 1).Value = "reMap"
$ws2.Cells.Item(21, 2).Value = 5.1
$ws2.Cells.Item(21, 3).Value = 3.92
$ws2.Cells.Item(21, 4).Value = 3.241
$ws2.Cells.Item(21, 5).Value = 3.2373
$ws2.Cells.Item(21, 6).Value = 3.23745
$ws2.Cells.Item(21, 7).Value = 3.23717
$ws2.Cells.Item(21, 8).Value = 3.2374
